$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column R (2021 data) mirroring the style of column Q for each row
$ws.Range("R4").Value = 2021
$ws.Range("R5").Value = 5.8
$ws.Range("R6").Value = 4.7
$ws.Range("R7").Value = 1.6
$ws.Range("R8").Value = 12.9
$ws.Range("R9").Value = 10.199999999999999
$ws.Range("R10").Value = 4.2
$ws.Range("R11").Value = 3.3
$ws.Range("R12").Value = 15.2
$ws.Range("R13").Value = 2.4
$ws.Range("R14").Value = 0.6

# Copy styles from column Q to column R for rows 4-14
for ($r = 4; $r -le 14; $r++) {
    $src = $ws.Range("Q$r")
    $dst = $ws.Range("R$r")
    $src.Copy()
    $dst.PasteSpecial(-4122) # xlPasteFormats
}

# Update the selection to match the target state
$ws.Range("T9").Select()
